$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emissions by Facility and Fuel")

for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 4).Value = "Tons CO2 eq"
}
